# Upgrade left table until javakheti:
# Add a "2023" column (K) to the Abasha sheet, mirroring the formatting of
# the existing last data column (J), and fill in the new year's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number format, borders, fonts, alignment, etc.) from
# the last existing year column (J3:J6) onto the new column (K3:K6) so the
# new cells look exactly like the rest of the table.
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new year header and data values.
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 569.5
$ws.Range("K5").Value = 351.4
$ws.Range("K6").Value = 640.79999999999995
